$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are plain text in the source data (e.g. "43.935.98" or
# "116.16") even when they look numeric, so force text formatting before
# assigning them to avoid Excel auto-converting to a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.935.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.315.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +22.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.38%  "
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.664.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.862"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.312.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.828.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000111"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +17.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.73%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "42.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.79%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "178.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0933"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.60%  "
$ws.Range("E37").Value = "  +3.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +22.37%  "
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.245"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.00%  "
$ws.Range("E41").Value = "  +2.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +14.54%  "
$ws.Range("B43").Value = "THORChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +15.87%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.81%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("E46").Value = "  +4.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.90%  "
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.48%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.467"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.35%  "
